{"js": "// Update the date and the 25 division problems in the table to the new\n// values, preserving all existing run/paragraph formatting.\nconst replacements = [\n  [\"2025-03-11 Tuesday\", \"2025-03-12 Wednesday\"],\n  [\"689\u00f74=\", \"962\u00f76=\"],\n  [\"446\u00f79=\", \"638\u00f79=\"],\n  [\"525\u00f77=\", \"352\u00f74=\"],\n  [\"575\u00f76=\", \"318\u00f75=\"],\n  [\"995\u00f74=\", \"682\u00f72=\"],\n  [\"744\u00f74=\", \"560\u00f76=\"],\n  [\"226\u00f79=\", \"454\u00f75=\"],\n  [\"157\u00f72=\", \"935\u00f72=\"],\n  [\"761\u00f76=\", \"243\u00f76=\"],\n  [\"469\u00f78=\", \"232\u00f79=\"],\n  [\"472\u00f76=\", \"443\u00f78=\"],\n  [\"688\u00f75=\", \"188\u00f79=\"],\n  [\"884\u00f73=\", \"995\u00f72=\"],\n  [\"910\u00f74=\", \"346\u00f72=\"],\n  [\"619\u00f73=\", \"252\u00f73=\"],\n  [\"115\u00f72=\", \"151\u00f74=\"],\n  [\"843\u00f72=\", \"217\u00f78=\"],\n  [\"124\u00f75=\", \"154\u00f79=\"],\n  [\"153\u00f72=\", \"920\u00f75=\"],\n  [\"115\u00f75=\", \"969\u00f73=\"],\n  [\"893\u00f76=\", \"485\u00f77=\"],\n  [\"975\u00f73=\", \"835\u00f75=\"],\n  [\"705\u00f79=\", \"161\u00f74=\"],\n  [\"379\u00f73=\", \"971\u00f79=\"],\n  [\"515\u00f77=\", \"345\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and the 25 division problems in the table to the new\n# values, preserving all existing run/paragraph formatting.\n$d = $word.ActiveDocument\n\n# NOTE: each pair is prefixed with the unary comma operator so PowerShell\n# keeps this as an array-of-arrays instead of flattening it into one\n# big array of strings.\n$pairs = @(\n    ,@(\"2025-03-11 Tuesday\", \"2025-03-12 Wednesday\")\n    ,@(\"689\u00f74=\", \"962\u00f76=\")\n    ,@(\"446\u00f79=\", \"638\u00f79=\")\n    ,@(\"525\u00f77=\", \"352\u00f74=\")\n    ,@(\"575\u00f76=\", \"318\u00f75=\")\n    ,@(\"995\u00f74=\", \"682\u00f72=\")\n    ,@(\"744\u00f74=\", \"560\u00f76=\")\n    ,@(\"226\u00f79=\", \"454\u00f75=\")\n    ,@(\"157\u00f72=\", \"935\u00f72=\")\n    ,@(\"761\u00f76=\", \"243\u00f76=\")\n    ,@(\"469\u00f78=\", \"232\u00f79=\")\n    ,@(\"472\u00f76=\", \"443\u00f78=\")\n    ,@(\"688\u00f75=\", \"188\u00f79=\")\n    ,@(\"884\u00f73=\", \"995\u00f72=\")\n    ,@(\"910\u00f74=\", \"346\u00f72=\")\n    ,@(\"619\u00f73=\", \"252\u00f73=\")\n    ,@(\"115\u00f72=\", \"151\u00f74=\")\n    ,@(\"843\u00f72=\", \"217\u00f78=\")\n    ,@(\"124\u00f75=\", \"154\u00f79=\")\n    ,@(\"153\u00f72=\", \"920\u00f75=\")\n    ,@(\"115\u00f75=\", \"969\u00f73=\")\n    ,@(\"893\u00f76=\", \"485\u00f77=\")\n    ,@(\"975\u00f73=\", \"835\u00f75=\")\n    ,@(\"705\u00f79=\", \"161\u00f74=\")\n    ,@(\"379\u00f73=\", \"971\u00f79=\")\n    ,@(\"515\u00f77=\", \"345\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n"}
